# "Completed the current dataset to use for the project. Converted images
# to BW, and added the testing set." -- appends the new testing-set rows
# (Image ID 124-160, worksheet rows 125-161) onto the bottom of Sheet1's
# Image-ID/Unit-Name table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is Image ID (col A) / Unit Name (col B) for the new rows.
$newRows = @(
    124, 'Scout'
    125, 'Scout'
    126, 'Scout'
    127, 'Scout'
    128, 'Intercessor'
    129, 'Intercessor'
    130, 'Intercessor'
    131, 'Intercessor'
    132, 'Librarian'
    133, 'Librarian'
    134, 'Eradicator'
    135, 'Eradicator'
    136, 'Eradicator'
    137, 'Eradicator'
    138, 'Bladeguard Veteran'
    139, 'Bladeguard Veteran'
    140, 'Bladeguard Veteran'
    141, 'Bladeguard Veteran'
    142, 'Bladeguard Veteran'
    143, 'Bladeguard Veteran'
    144, 'Bladeguard Veteran'
    145, 'Gladiator Lancer'
    146, 'Gladiator Lancer'
    147, 'Impulsor'
    148, 'Impulsor'
    149, 'Eliminator'
    150, 'Eliminator'
    151, 'Reiver'
    152, 'Reiver'
    153, 'Reiver'
    154, 'Reiver'
    155, 'Captain In Jump Pack'
    156, 'Captain In Jump Pack'
    157, 'Sternguard Veteran'
    158, 'Sternguard Veteran'
    159, 'Storm Speeder Thunderstrike'
    160, 'Storm Speeder Thunderstrike'
)

$startRow = 125
for ($i = 0; $i -lt $newRows.Count; $i += 2) {
    $r = $startRow + ($i / 2)
    $ws.Range("A$r").Value = $newRows[$i]
    $ws.Range("B$r").Value = $newRows[$i + 1]
}

# Move the view / selection to reflect where the new data ends, same as
# the author's saved window state (topLeftCell A142, active cell C159).
$excel.ActiveWindow.ScrollRow = 142
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C159").Select()
